# Update with GW7 data
# Applies new 2023-24 season values (GW7 instead of GW6) to the Midfielders
# analysis sheet, and recolors the dependent "traffic-light" columns
# (F: mean total vs. effect-size sign, K: BH-significance flag,
#  L: effect-size-label, N: average opposition strength heat scale)
# to match the colors implied by the new data.

function Get-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Color palette (kept identical to the values already present in the
# workbook's style table; only which rows use which color changes).
# ---------------------------------------------------------------------------
$White  = Get-RGB 0xFF 0xFF 0xFF
$Black  = Get-RGB 0x00 0x00 0x00
$OffWht = Get-RGB 0xF1 0xF1 0xF1

$Green  = Get-RGB 0x22 0x8B 0x22   # 00228B22 - positive effect size (F)
$Red    = Get-RGB 0xDC 0x14 0x3C   # 00DC143C - negative effect size (F) / not-significant (K)

$LblVerySmall = Get-RGB 0xCC 0xEC 0xE6   # 00CCECE6
$LblSmall     = Get-RGB 0x98 0xD8 0xC9   # 0098D8C9
$LblMedium    = Get-RGB 0x65 0xC2 0xA3   # 0065C2A3
$LblLarge     = Get-RGB 0x40 0xAD 0x75   # 0040AD75

$OppLow   = Get-RGB 0x3A 0x85 0x5E   # 003A855E - 2.857142857142857
$OppMid1  = Get-RGB 0xE0 0xEB 0xE4   # 00E0EBE4 - 3
$OppMid2  = Get-RGB 0xDE 0x56 0x60   # 00DE5660 - 3.142857142857143
$OppHigh  = Get-RGB 0xDA 0x3B 0x46   # 00DA3B46 - 3.166666666666667

# Label -> (fill, font) lookup for column L
$LabelFill = @{
    "Very small" = $LblVerySmall
    "Small"      = $LblSmall
    "Medium"     = $LblMedium
    "Large"      = $LblLarge
}

# Opposition-strength value -> fill color lookup for column N
$OppFill = @{
    "2.857142857142857" = $OppLow
    "3"                  = $OppMid1
    "3.142857142857143" = $OppMid2
    "3.166666666666667" = $OppHigh
}
$OppFont = @{
    "2.857142857142857" = $OffWht
    "3"                  = $Black
    "3.142857142857143" = $OffWht
    "3.166666666666667" = $OffWht
}

# ---------------------------------------------------------------------------
# New GW7 data, one row per player (row numbers match the worksheet).
# Columns: E,F,G,H,I,J,K,L,M(optional),N
# ---------------------------------------------------------------------------
$rows = @(
    @{R=2;  E=7; F=5.571428571428571;  G=-0.05536909615918346; H=1.60819424531872;  I=-0.146492858755621;  J=0.4441652170599477;  K=$false; L="Very small"; N=3.142857142857143},
    @{R=3;  E=7; F=2.285714285714286;  G=-0.7555904464118185;  H=1.226112954436575; I=-1.999104414221948;   J=0.04627054157368334; K=$false; L="Medium";     M=7.7;  N=3.142857142857143},
    @{R=4;  E=7; F=3.714285714285714;  G=-0.1238191175791926;  H=0.9524884786553296;I=-0.3275945926700094;  J=0.3771734669037085;  K=$false; L="Very small"; N=3},
    @{R=5;  E=7; F=6;                  G=0.6312062720519506;   H=1.223857191206658; I=1.670014821833641;    J=0.0729771728275128;  K=$false; L="Medium";     M=7.2;  N=3.142857142857143},
    @{R=6;  E=7; F=6.428571428571429;  G=0.7785568269322275;   H=1.204544812561751; I=2.059867745594229;    J=0.04253601048153195; K=$false; L="Medium";     N=3.142857142857143},
    @{R=7;  E=7; F=0.1428571428571428; G=-0.9115936300930927;  H=1.54832836810122;  I=-2.41185004197693;    J=0.026220811648435;   K=$false; L="Large";      N=2.857142857142857},
    @{R=8;  E=7; F=6.571428571428571;  G=0.914613230211084;    H=1.149644547171203; I=2.419839152947996;    J=0.0259372357625665;  K=$false; L="Large";      M=9.199999999999999; N=3.142857142857143},
    @{R=9;  E=7; F=3.428571428571428;  G=-0.1095929248606583;  H=1.153923006626749; I=-0.28995562463349;    J=0.3908049826392022;  K=$false; L="Very small"; N=3.142857142857143},
    @{R=10; E=7; F=5.285714285714286;  G=0.5440810031322788;   H=1.12472678992139;  I=1.439503027362564;    J=0.1000342397129754;  K=$false; L="Medium";     N=3.142857142857143},
    @{R=11; E=7; F=3.714285714285714;  G=0.03010563150430234;  H=1.258607554324357; I=0.07965201402293536;  J=0.4695521183487202;  K=$false; L="Very small"; N=3},
    @{R=12; E=7; F=5.142857142857143;  G=0.453263759121011;    H=1.333361057222804; I=1.19922318495248;     J=0.1378236780625797;  K=$false; L="Small";      N=3.142857142857143},
    @{R=13; E=7; F=5.166666666666667;  G=0.6090531025646843;   H=1.023408024780048; I=1.61140304461847;     J=0.07910904855933919; K=$false; L="Medium";     M=6.3;  N=3.166666666666667},
    @{R=14; E=7; F=4.857142857142857;  G=0.4823067903233182;   H=1.094472326824878; I=1.276063822833274;    J=0.1245509421821516;  K=$false; L="Small";      N=3.142857142857143},
    @{R=15; E=7; F=3.857142857142857;  G=0.1279676139908083;   H=1.210307623619678; I=0.3385704824899884;   J=0.373233581234662;   K=$false; L="Very small"; M=6.4;  N=3.142857142857143},
    @{R=16; E=7; F=1.428571428571429;  G=-0.6302907132720854;  H=1.184304538351127; I=-1.667592480991456;   J=0.07322128002496231; K=$false; L="Medium";     N=3.142857142857143},
    @{R=17; E=7; F=3.714285714285714;  G=0.09789371633381265;  H=1.335372273832734; I=0.25900242833517;     J=0.4021450606515578;  K=$false; L="Very small"; M=7.1;  N=3},
    @{R=18; E=7; F=4.857142857142857;  G=0.5298144563823256;   H=1.080813063708611; I=1.401757292594511;    J=0.1052708972960017;  K=$false; L="Medium";     M=6.6;  N=3.142857142857143},
    @{R=19; E=7; F=2.142857142857143;  G=-0.3650440845146808;  H=1.214447156117957; I=-0.96581586520109;    J=0.1857175458037973;  K=$false; L="Small";      M=6.3;  N=3.142857142857143},
    @{R=20; E=7; F=4.285714285714286;  G=0.2855900012824901;   H=1.364912859417114; I=0.7556001203200863;   J=0.2392461003322635;  K=$false; L="Small";      N=2.857142857142857},
    @{R=21; E=7; F=3.571428571428572;  G=0.1046866837448216;   H=1.160496982655813; I=0.2769749307688659;   J=0.3955470658989948;  K=$false; L="Very small"; N=3.142857142857143}
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    if ($row.ContainsKey("M")) {
        $ws.Range("M$r").Value = $row.M
    }
    $ws.Range("N$r").Value = $row.N

    # --- recolor F: green when effect size (G) is non-negative, red otherwise ---
    $fCell = $ws.Range("F$r")
    if ($row.G -ge 0) {
        $fCell.Interior.Color = $Green
    } else {
        $fCell.Interior.Color = $Red
    }
    $fCell.Font.Color = $White

    # --- recolor K: always red / white text (BH significance is False for all rows) ---
    $kCell = $ws.Range("K$r")
    if ($row.K) {
        $kCell.Interior.Color = $Green
    } else {
        $kCell.Interior.Color = $Red
    }
    $kCell.Font.Color = $White

    # --- recolor L based on the effect-size label ---
    $lCell = $ws.Range("L$r")
    $lCell.Interior.Color = $LabelFill[$row.L]
    $lCell.Font.Color = $White

    # --- recolor N based on the average-opposition-strength heat scale ---
    $nKey = [string]$row.N
    $nCell = $ws.Range("N$r")
    $nCell.Interior.Color = $OppFill[$nKey]
    $nCell.Font.Color = $OppFont[$nKey]
}
